{"js": "// This script updates the 25 \"two-digit \u00f7 one-digit\" practice answers in the\n// single table on the page. The table has 20 rows (5 data rows holding the\n// problems, interleaved with 3 blank spacer rows), 5 columns each \u2014 so the\n// 25 populated cells are addressed directly by (row, col) to avoid any\n// ambiguity from duplicate/overlapping text values across cells.\nconst cellEdits = [\n  { row: 0, col: 0, oldText: \"62\u00f78=7, 6\", newText: \"55\u00f74=13, 3\" },\n  { row: 0, col: 1, oldText: \"68\u00f73=22, 2\", newText: \"67\u00f75=13, 2\" },\n  { row: 0, col: 2, oldText: \"30\u00f75=6, 0\", newText: \"28\u00f77=4, 0\" },\n  { row: 0, col: 3, oldText: \"46\u00f76=7, 4\", newText: \"85\u00f78=10, 5\" },\n  { row: 0, col: 4, oldText: \"58\u00f75=11, 3\", newText: \"96\u00f73=32, 0\" },\n  { row: 4, col: 0, oldText: \"90\u00f74=22, 2\", newText: \"27\u00f72=13, 1\" },\n  { row: 4, col: 1, oldText: \"51\u00f72=25, 1\", newText: \"12\u00f77=1, 5\" },\n  { row: 4, col: 2, oldText: \"73\u00f79=8, 1\", newText: \"82\u00f72=41, 0\" },\n  { row: 4, col: 3, oldText: \"45\u00f73=15, 0\", newText: \"93\u00f77=13, 2\" },\n  { row: 4, col: 4, oldText: \"28\u00f75=5, 3\", newText: \"36\u00f78=4, 4\" },\n  { row: 8, col: 0, oldText: \"63\u00f78=7, 7\", newText: \"63\u00f77=9, 0\" },\n  { row: 8, col: 1, oldText: \"56\u00f76=9, 2\", newText: \"41\u00f73=13, 2\" },\n  { row: 8, col: 2, oldText: \"46\u00f77=6, 4\", newText: \"70\u00f72=35, 0\" },\n  { row: 8, col: 3, oldText: \"50\u00f79=5, 5\", newText: \"60\u00f77=8, 4\" },\n  { row: 8, col: 4, oldText: \"45\u00f75=9, 0\", newText: \"14\u00f77=2, 0\" },\n  { row: 12, col: 0, oldText: \"91\u00f73=30, 1\", newText: \"61\u00f72=30, 1\" },\n  { row: 12, col: 1, oldText: \"83\u00f79=9, 2\", newText: \"60\u00f77=8, 4\" },\n  { row: 12, col: 2, oldText: \"79\u00f77=11, 2\", newText: \"28\u00f76=4, 4\" },\n  { row: 12, col: 3, oldText: \"68\u00f76=11, 2\", newText: \"95\u00f75=19, 0\" },\n  { row: 12, col: 4, oldText: \"24\u00f75=4, 4\", newText: \"45\u00f74=11, 1\" },\n  { row: 16, col: 0, oldText: \"39\u00f74=9, 3\", newText: \"57\u00f78=7, 1\" },\n  { row: 16, col: 1, oldText: \"36\u00f78=4, 4\", newText: \"79\u00f76=13, 1\" },\n  { row: 16, col: 2, oldText: \"63\u00f72=31, 1\", newText: \"78\u00f73=26, 0\" },\n  { row: 16, col: 3, oldText: \"44\u00f77=6, 2\", newText: \"50\u00f77=7, 1\" },\n  { row: 16, col: 4, oldText: \"85\u00f79=9, 4\", newText: \"87\u00f73=29, 0\" },\n];\n\n\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Resolve every target paragraph range up front (load each cell's text so we\n// can sanity-check against the expected \"before\" value before mutating it).\nconst targets = cellEdits.map((edit) => {\n  const cell = table.getCell(edit.row, edit.col);\n  const para = cell.body.paragraphs.getFirst();\n  para.load(\"text\");\n  return { edit, para };\n});\nawait context.sync();\n\nfor (const { edit, para } of targets) {\n  if (para.text !== edit.oldText) {\n    throw new Error(\n      `Unexpected cell text at row ${edit.row}, col ${edit.col}: ` +\n      `expected \"${edit.oldText}\" but found \"${para.text}\"`\n    );\n  }\n}\n\n// Replace each paragraph's text in place via its Range so the existing run\n// formatting (font, size) and paragraph formatting (alignment) are kept\n// untouched \u2014 only the visible digits/numbers change, matching the diff.\nfor (const { para, edit } of targets) {\n  para.getRange().insertText(edit.newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# (row, col, expected \"before\" text, new \"after\" text) for the 25 populated\n# problem cells -- the table has 20 rows (5 data rows + 3 blank spacer rows\n# after each), 5 columns, addressed directly by position so there is no\n# ambiguity even though a couple of the new values collide with other cells' old values.\n$cellEdits = @(\n    @{ Row = 1; Col = 1; OldText = \"62\u00f78=7, 6\"; NewText = \"55\u00f74=13, 3\" }\n    @{ Row = 1; Col = 2; OldText = \"68\u00f73=22, 2\"; NewText = \"67\u00f75=13, 2\" }\n    @{ Row = 1; Col = 3; OldText = \"30\u00f75=6, 0\"; NewText = \"28\u00f77=4, 0\" }\n    @{ Row = 1; Col = 4; OldText = \"46\u00f76=7, 4\"; NewText = \"85\u00f78=10, 5\" }\n    @{ Row = 1; Col = 5; OldText = \"58\u00f75=11, 3\"; NewText = \"96\u00f73=32, 0\" }\n    @{ Row = 5; Col = 1; OldText = \"90\u00f74=22, 2\"; NewText = \"27\u00f72=13, 1\" }\n    @{ Row = 5; Col = 2; OldText = \"51\u00f72=25, 1\"; NewText = \"12\u00f77=1, 5\" }\n    @{ Row = 5; Col = 3; OldText = \"73\u00f79=8, 1\"; NewText = \"82\u00f72=41, 0\" }\n    @{ Row = 5; Col = 4; OldText = \"45\u00f73=15, 0\"; NewText = \"93\u00f77=13, 2\" }\n    @{ Row = 5; Col = 5; OldText = \"28\u00f75=5, 3\"; NewText = \"36\u00f78=4, 4\" }\n    @{ Row = 9; Col = 1; OldText = \"63\u00f78=7, 7\"; NewText = \"63\u00f77=9, 0\" }\n    @{ Row = 9; Col = 2; OldText = \"56\u00f76=9, 2\"; NewText = \"41\u00f73=13, 2\" }\n    @{ Row = 9; Col = 3; OldText = \"46\u00f77=6, 4\"; NewText = \"70\u00f72=35, 0\" }\n    @{ Row = 9; Col = 4; OldText = \"50\u00f79=5, 5\"; NewText = \"60\u00f77=8, 4\" }\n    @{ Row = 9; Col = 5; OldText = \"45\u00f75=9, 0\"; NewText = \"14\u00f77=2, 0\" }\n    @{ Row = 13; Col = 1; OldText = \"91\u00f73=30, 1\"; NewText = \"61\u00f72=30, 1\" }\n    @{ Row = 13; Col = 2; OldText = \"83\u00f79=9, 2\"; NewText = \"60\u00f77=8, 4\" }\n    @{ Row = 13; Col = 3; OldText = \"79\u00f77=11, 2\"; NewText = \"28\u00f76=4, 4\" }\n    @{ Row = 13; Col = 4; OldText = \"68\u00f76=11, 2\"; NewText = \"95\u00f75=19, 0\" }\n    @{ Row = 13; Col = 5; OldText = \"24\u00f75=4, 4\"; NewText = \"45\u00f74=11, 1\" }\n    @{ Row = 17; Col = 1; OldText = \"39\u00f74=9, 3\"; NewText = \"57\u00f78=7, 1\" }\n    @{ Row = 17; Col = 2; OldText = \"36\u00f78=4, 4\"; NewText = \"79\u00f76=13, 1\" }\n    @{ Row = 17; Col = 3; OldText = \"63\u00f72=31, 1\"; NewText = \"78\u00f73=26, 0\" }\n    @{ Row = 17; Col = 4; OldText = \"44\u00f77=6, 2\"; NewText = \"50\u00f77=7, 1\" }\n    @{ Row = 17; Col = 5; OldText = \"85\u00f79=9, 4\"; NewText = \"87\u00f73=29, 0\" }\n)\n\nforeach ($edit in $cellEdits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $edit.OldText) {\n        throw \"Unexpected cell text at row $($edit.Row), col $($edit.Col): expected `\"$($edit.OldText)`\" but found `\"$current`\"\"\n    }\n}\n\nforeach ($edit in $cellEdits) {\n    $cell = $t.Cell($edit.Row, $edit.Col)\n    $cell.Range.Text = $edit.NewText\n}\n"}
